$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("利欧股份", "利欧股份", "利欧股份"),
    @("金风科技", "航天电子", "金风科技"),
    @("蓝色光标", "金风科技", "岩山科技"),
    @("航天电子", "蓝色光标", "航天发展"),
    @("岩山科技", "岩山科技", "航天电子"),
    @("易点天下", "昆仑万维", "蓝色光标"),
    @("中国卫通", "中国卫通", "海格通信"),
    @("航天发展", "乾照光电", "中国卫通"),
    @("通宇通讯", "易点天下", "雷科防务"),
    @("雷科防务", "省广集团", "通宇通讯"),
    @("乾照光电", "雷科防务", "中国卫星"),
    @("引力传媒", "航天发展", "银河电子"),
    @("华胜天成", "通宇通讯", "雪人集团"),
    @("中国卫星", "中航光电", "万向钱潮"),
    @("昆仑万维", "华胜天成", "华胜天成"),
    @("雪人集团", "美年健康", "巨力索具"),
    @("银河电子", "安泰科技", "永鼎股份"),
    @("巨力索具", "中国卫星", "乾照光电"),
    @("省广集团", "中国一重", "天普股份"),
    @("鲁信创投", "巨力索具", "天下秀")
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
